# Applies the "Added Heuristic to report" edit:
#  1. Marks the three inline-picture runs as NoProofing (adds <w:rPr><w:noProof/></w:rPr>).
#  2. Splits the trailing (previously empty) Heading1 paragraph into:
#       - a blank paragraph indented like the surrounding body text,
#       - a centred "Heuristic" heading,
#       - a new body paragraph explaining the A* heuristic, with the
#         _GoBack bookmark relocated into the middle of that paragraph.

$d = $word.ActiveDocument

# --- 1. Inline pictures: turn on NoProofing for the run that hosts each <w:drawing> ---
for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
    $shape = $d.InlineShapes.Item($i)
    $shape.Range.NoProofing = -1
}

# --- 2. Find the empty "Heading 1" paragraph that currently just holds the _GoBack bookmark ---
$headingPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Style.NameLocal -eq "Heading 1" -and $candidate.Range.Text.Trim() -eq "") {
        $headingPara = $candidate
    }
}

# Add a blank paragraph in front of it; this becomes the plain indented <w:p> that now
# separates the last body paragraph from the new "Heuristic" heading.
$headingPara.Range.InsertParagraphBefore()

$headingIdx = -1
$blankIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Style.NameLocal -eq "Heading 1") {
        $headingIdx = $i
        $blankIdx = $i - 1
    }
}

$blankPara = $d.Paragraphs.Item($blankIdx)
$blankPara.Style = "Normal"
$blankPara.LeftIndent = 18

# The _GoBack bookmark used to sit in the (now blank) heading paragraph; drop it here, it gets
# re-created further down, in the middle of the new justification paragraph.
$d.Bookmarks.Item("_GoBack").Delete()

$headingPara2 = $d.Paragraphs.Item($headingIdx)
$headingPara2.Alignment = 1
$headingPara2.Range.InsertBefore("Heuristic")

# --- 3. Turn the two trailing empty paragraphs into a single paragraph with the new text ---
$para1Idx = $headingIdx + 1
$para2Idx = $headingIdx + 2

$run1 = "As described in solve_sokoban_macro, we have utilised the A* graph search algorithm to find a list of steps for the worker to move through the warehouse with a crate to the goal. A* is an informed search algorithm, also known as a best-first search. A* starts from a specific starting node of a graph and aims to find a path to the given goal with the smallest cost. In our case, the smallest cost is calc"
$run2 = "ulated as the least distance travelled. This is done by "
$run3 = "maintain"
$run4 = "ing"
$run5 = " a tree of paths that originates at the starting node and extending those paths one edge at a time until it reaches its goal."
$fullText = $run1 + $run2 + $run3 + $run4 + $run5

$targetPara = $d.Paragraphs.Item($para2Idx)
$targetPara.Range.InsertBefore($fullText)

# Drop the now-redundant leading empty paragraph so the text above becomes a single paragraph.
$leadingPara = $d.Paragraphs.Item($para1Idx)
$leadingPara.Range.Delete()

$finalPara = $d.Paragraphs.Item($para1Idx)
$paraStart = $finalPara.Range.Start

# --- 4. Recreate the run boundaries seen in the target (matching separately-typed runs) and
#        put the _GoBack bookmark back, in the middle of the paragraph this time. A temporary
#        bookmark forces Word to break the run at that character offset; deleting the temporary
#        bookmark again leaves the run split in place without leaving stray bookmarks behind. ---
function Split-RunAt($pos) {
    $tmpRange = $d.Range($pos, $pos)
    $d.Bookmarks.Add("TmpRunSplit", $tmpRange) | Out-Null
    $d.Bookmarks.Item("TmpRunSplit").Delete()
}

$b1 = $paraStart + $run1.Length
$b2 = $b1 + $run2.Length
$b3 = $b2 + $run3.Length
$b4 = $b3 + $run4.Length

Split-RunAt $b2
Split-RunAt $b3
Split-RunAt $b4

$bookmarkRange = $d.Range($b1, $b1)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null
